$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain Text so numeric-looking values (e.g. "0.999", "2.09")
# are not auto-converted to numbers by Excel when assigned below.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '37.851.08'
$ws.Range("E2").Value = '  +1.57%  '
$ws.Range("D3").Value = '2.048.48'
$ws.Range("E3").Value = '  +0.98%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '229.48'
$ws.Range("E5").Value = '  +0.86%  '
$ws.Range("D6").Value = '0.615'
$ws.Range("E6").Value = '  +1.97%  '
$ws.Range("D7").Value = '58.07'
$ws.Range("E7").Value = '  +4.91%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  +1.77%  '
$ws.Range("E10").Value = '  +2.24%  '
$ws.Range("E11").Value = '  +0.39%  '
$ws.Range("D12").Value = '2.353.36'
$ws.Range("E12").Value = '  +1.14%  '
$ws.Range("D13").Value = '14.54'
$ws.Range("E13").Value = '  +1.76%  '
$ws.Range("D14").Value = '20.66'
$ws.Range("E14").Value = '  +1.30%  '
$ws.Range("D15").Value = '5.27'
$ws.Range("E15").Value = '  +1.87%  '
$ws.Range("D16").Value = '0.748'
$ws.Range("E16").Value = '  +0.44%  '
$ws.Range("D17").Value = '2.055.38'
$ws.Range("E17").Value = '  +1.59%  '
$ws.Range("D18").Value = '37.821.84'
$ws.Range("E18").Value = '  +1.61%  '
$ws.Range("E19").Value = '  -2.41%  '
$ws.Range("D20").Value = '69.58'
$ws.Range("E20").Value = '  +0.44%  '
$ws.Range("E21").Value = '  +1.23%  '
$ws.Range("D22").Value = '224.48'
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("E24").Value = '  +0.33%  '
$ws.Range("D25").Value = '2.24'
$ws.Range("E25").Value = '  +1.83%  '
$ws.Range("D26").Value = '166.68'
$ws.Range("E26").Value = '  +0.72%  '
$ws.Range("D27").Value = '9.25'
$ws.Range("E27").Value = '  -0.77%  '
$ws.Range("D28").Value = '0.132'
$ws.Range("E28").Value = '  +3.64%  '
$ws.Range("D29").Value = '19.00'
$ws.Range("E29").Value = '  +0.89%  '
$ws.Range("E30").Value = '  -0.97%  '
$ws.Range("E31").Value = '  +1.37%  '
$ws.Range("D32").Value = '4.53'
$ws.Range("E32").Value = '  -0.14%  '
$ws.Range("D33").Value = '2.09'
$ws.Range("E33").Value = '  +13.47%  '
$ws.Range("D34").Value = '4.57'
$ws.Range("E34").Value = '  +2.37%  '
$ws.Range("D35").Value = '0.0611'
$ws.Range("E35").Value = '  -1.16%  '
$ws.Range("E36").Value = '  -1.55%  '
$ws.Range("D37").Value = '5.96'
$ws.Range("E37").Value = '  +8.93%  '
$ws.Range("D38").Value = '3.31'
$ws.Range("E38").Value = '  +4.36%  '
$ws.Range("E39").Value = '  +0.12%  '
$ws.Range("E40").Value = '  +0.33%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '1.483.21'
$ws.Range("E41").Value = '  +0.14%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '97.49'
$ws.Range("E42").Value = '  +1.33%  '
$ws.Range("D43").Value = '2.85'
$ws.Range("E43").Value = '  +2.88%  '
$ws.Range("D44").Value = '0.0933'
$ws.Range("E44").Value = '  +1.04%  '
$ws.Range("D45").Value = '16.56'
$ws.Range("E45").Value = '  +0.67%  '
$ws.Range("D46").Value = '4.21'
$ws.Range("E46").Value = '  +16.05%  '
$ws.Range("E47").Value = '  -0.50%  '
$ws.Range("D48").Value = '1.01'
$ws.Range("E48").Value = '  -0.59%  '
$ws.Range("E49").Value = '  +1.31%  '
$ws.Range("D50").Value = '6.98'
$ws.Range("E50").Value = '  -3.33%  '
$ws.Range("D51").Value = '2.242.99'
$ws.Range("E51").Value = '  +1.34%  '
